$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-26 Sunday", "2024-05-27 Monday"),
    @("91×67=6097", "41×67=2747"),
    @("56×20=1120", "38×62=2356"),
    @("80×73=5840", "71×70=4970"),
    @("97×18=1746", "81×83=6723"),
    @("34×51=1734", "71×34=2414"),
    @("13×40=520", "16×89=1424"),
    @("76×89=6764", "51×73=3723"),
    @("56×86=4816", "23×65=1495"),
    @("27×27=729", "94×31=2914"),
    @("16×29=464", "42×53=2226"),
    @("59×48=2832", "37×76=2812"),
    @("51×14=714", "74×28=2072"),
    @("74×25=1850", "23×70=1610"),
    @("56×75=4200", "66×95=6270"),
    @("95×16=1520", "89×95=8455"),
    @("31×53=1643", "61×63=3843"),
    @("14×69=966", "44×44=1936"),
    @("99×97=9603", "64×67=4288"),
    @("51×62=3162", "81×68=5508"),
    @("50×85=4250", "68×49=3332"),
    @("76×14=1064", "18×89=1602"),
    @("54×97=5238", "27×82=2214"),
    @("37×38=1406", "11×91=1001"),
    @("60×23=1380", "69×44=3036"),
    @("12×20=240", "71×33=2343")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
